$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s0 = '[Jianghong%An%NULL%1,               Xuejiao%Liao%NULL%1,               Tongyang%Xiao%NULL%1,               Shen%Qian%NULL%1,               Jing%Yuan%NULL%1,               Haocheng%Ye%NULL%1,               Furong%Qi%NULL%1,               Chengguang%Shen%NULL%1,               Lifei%Wang%NULL%1,               Yang%Liu%NULL%1,               Xiaoya%Cheng%NULL%1,               Na%Li%NULL%1,               Qingxian%Cai%NULL%1,               Fang%Wang%NULL%1,               Jun%Chen%NULL%0,               Guojun%Li%NULL%1,               Qiu’e%Cai%NULL%1,               Yingxia%Liu%NULL%0,               Yunfang%Wang%NULL%1,               Feng%Zhang%NULL%1,               Yang%Fu%NULL%1,               Qing%He%NULL%1,               Xiaohua%Tan%NULL%1,               Lei%Liu%NULL%0,               Zheng%Zhang%NULL%1]'
$s1 = '[Jie%Chen%NULL%2,               Xiaoping%Xu%NULL%1,               Jing%Hu%NULL%1,               Qiangda%Chen%NULL%1,               Fengfeng%Xu%NULL%1,               Hui%Liang%NULL%1,               Nanmei%Liu%NULL%1,               Hengmei%Zhu%NULL%1,               Jinlong%Lan%NULL%1,               Lan%Zhou%NULL%1,               Jiajun%Xing%NULL%1,               Ning%Pu%NULL%1,               Zhigang%Cai%NULL%1]'
$s2 = '[Q. Sue%Huang%NULL%1,               Tim%Wood%NULL%1,               Lauren%Jelley%NULL%1,               Tineke%Jennings%NULL%1,               Sarah%Jefferies%NULL%1,               Karen%Daniells%NULL%1,               Annette%Nesdale%NULL%1,               Tony%Dowell%NULL%1,               Nikki%Turner%NULL%1,               Priscilla%Campbell-Stokes%NULL%1,               Michelle%Balm%NULL%1,               Hazel C%Dobinson%NULL%1,               Cameron C.%Grant%NULL%1,               Shelley%James%NULL%1,               Nayyereh%Aminisani%NULL%1,               Jacqui%Ralston%NULL%1,               Wendy%Gunn%NULL%1,               Judy%Bocacao%NULL%1,               Jessica%Danielewicz%NULL%1,               Tessa%Moncrieff%NULL%1,               Andrea%McNeill%NULL%1,               Liza%Lopez%NULL%1,               Ben%Waite%NULL%1,               Tomasz%Kiedrzynski%NULL%1,               Hannah%Schrader%NULL%1,               Rebekah%Gray%NULL%1,               Kayla%Cook%NULL%1,               Danielle%Currin%NULL%1,               Chaune%Engelbrecht%NULL%1,               Whitney%Tapurau%NULL%1,               Leigh%Emmerton%NULL%1,               Maxine%Martin%NULL%1,               Michael G.%Baker%NULL%1,               Susan%Taylor%NULL%1,               Adrian%Trenholme%NULL%1,               Conroy%Wong%NULL%1,               Shirley%Lawrence%NULL%1,               Colin%McArthur%NULL%1,               Alicia%Stanley%NULL%1,               Sally%Roberts%NULL%1,               Fahimeh%Ranama%NULL%1,               Jenny%Bennett%NULL%1,               Chris%Mansell%NULL%1,               Meik%Dilcher%NULL%1,               Anja%Werno%NULL%1,               Jennifer%Grant%NULL%1,               Antje%van der Linden%NULL%1,               Ben%Youngblood%NULL%1,               Paul G.%Thomas%NULL%1,               Richard J.%Webby%NULL%1,                S.%Huang%null%1,                T.% Wood%null%1,                L.% Jelley%null%1,                T.% Jennings%null%1,                S.% Jeffries%null%1,                K.% Daniels%null%1,                A.% Nesdale%null%1,                T.% Dowell%null%1,                N.% Turner%null%1,                P.% Campbell-Stokes%null%1,                M.% Balm%null%1,                H. C.% Dobinson%null%1,                C. C.% Grant%null%1,                S.% James%null%1,                N.% Aminisani%null%1,                J.% Ralston%null%1,                W.% Gunn%null%1,                J.% Bucacao%null%1,                J.% Danielewicz%null%1,                T.% Moncrieff%null%1,                A.% McNeill%null%1,                L.% Lopez%null%1,                B.% Waite%null%1,                T.% Kiedrzynski%null%1,                H.% Schrader%null%1,                R.% Gray%null%1,                K.% Cook%null%1,                D.% Currin%null%1,                C.% Engelbrecht%null%1,                W.% Tapurau%null%1,                L.% Emmerton%null%1,                M.% Martin%null%1,                M. G.% Baker%null%1,                S.% Taylor%null%1,                A.% Trenholme%null%1,                C.% Wong%null%1,                S.% Lawrence%null%1,                C.% McArthur%null%1,                A.% Stanley%null%1,                S.% Roberts%null%1,                F.% Ranama%null%1,                J.% Bennett%null%1,                C.% Mansell%null%1,                M.% Dilcher%null%1,                A.% Werno%null%1,                J.% Grant%null%1,               couldn''t find firstname% v %null%1,              S.%Huang%null%1,              T.% Wood%null%1,              L.% Jelley%null%1,              T.% Jennings%null%1,              S.% Jeffries%null%1,              K.% Daniels%null%1,              A.% Nesdale%null%1,              T.% Dowell%null%1,              N.% Turner%null%1,              P.% Campbell-Stokes%null%1,              M.% Balm%null%1,              H. C.% Dobinson%null%1,              C. C.% Grant%null%1,              S.% James%null%1,              N.% Aminisani%null%1,              J.% Ralston%null%1,              W.% Gunn%null%1,              J.% Bucacao%null%1,              J.% Danielewicz%null%1,              T.% Moncrieff%null%1,              A.% McNeill%null%1,              L.% Lopez%null%1,              B.% Waite%null%1,              T.% Kiedrzynski%null%1,              H.% Schrader%null%1,              R.% Gray%null%1,              K.% Cook%null%1,              D.% Currin%null%1,              C.% Engelbrecht%null%1,              W.% Tapurau%null%1,              L.% Emmerton%null%1,              M.% Martin%null%1,              M. G.% Baker%null%1,              S.% Taylor%null%1,              A.% Trenholme%null%1,              C.% Wong%null%1,              S.% Lawrence%null%1,              C.% McArthur%null%1,              A.% Stanley%null%1,              S.% Roberts%null%1,              F.% Ranama%null%1,              J.% Bennett%null%1,              C.% Mansell%null%1,              M.% Dilcher%null%1,              A.% Werno%null%1,              J.% Grant%null%1,             couldn''t find firstname% v %null%1]'
$s3 = '[Minlin%Jiang%NULL%1,               Ya%Li%NULL%1,               Mingli%Han%NULL%1,               Zhenhua%Wang%NULL%1,               Yuhang%Zhang%NULL%1,               Xinwei%Du%NULL%1]'
$s4 = '[Hong%Zheng%NULL%1,               Yiran%Hu%NULL%1,               Li%Dong%NULL%1,               Qi%Shu%NULL%1,               Mingyang%Zhu%NULL%1,               Yuping%Li%NULL%1,               Chengshui%Chen%NULL%1,               Hongchang%Gao%gaohc27@wmu.edu.cn%1,               Li%Yang%taiyang2630@163.com%0,               Li%Yang%taiyang2630@163.com%0]'
$s5 = '[Youjiang%Li%NULL%2,               Yingying%Hu%NULL%2,               Yingying%Hu%NULL%0,               Yuanyuan%Yu%NULL%1,               Xiaodong%Zhang%NULL%1,               Bin%Li%NULL%1,               Jianguo%Wu%NULL%1,               Junyu%Li%NULL%1,               Yingping%Wu%NULL%1,               Xiaoping%Xia%NULL%2,               Xiaoping%Xia%NULL%0,               Huina%Tang%NULL%1,               Jian%Xu%xuj@zju.edu.cn%1]'
$s6 = '[Yun%Ling%NULL%0,               Shui-Bao%Xu%NULL%0,               Yi-Xiao%Lin%NULL%0,               Di%Tian%NULL%0,               Zhao-Qin%Zhu%NULL%0,               Fa-Hui%Dai%NULL%0,               Fan%Wu%NULL%0,               Zhi-Gang%Song%NULL%0,               Wei%Huang%NULL%0,               Jun%Chen%NULL%0,               Bi-Jie%Hu%NULL%0,               Sheng%Wang%NULL%0,               En-Qiang%Mao%NULL%0,               Lei%Zhu%NULL%0,               Wen-Hong%Zhang%NULL%0,               Hong-Zhou%Lu%NULL%0,               Peng%Lyu%NULL%0,               Peng%Lyu%NULL%0]'
$s7 = '[Ai Tang%Xiao%NULL%0,               Yi Xin%Tong%NULL%0,               Sheng%Zhang%aloof3737@126.com%0]'
$s8 = '[Justin%Wong%justin.wong@moh.gov.bn%1,               Wee Chian%Koh%NULL%2,               Wee Chian%Koh%NULL%0,               Riamiza Natalie%Momin%NULL%2,               Riamiza Natalie%Momin%NULL%0,               Mohammad Fathi%Alikhan%NULL%1,               Noraskhin%Fadillah%NULL%1,               Lin%Naing%NULL%1]'
$s9 = '[Jinru%Wu%NULL%1,               Xinyi%Liu%NULL%1,               Jianjun%Liu%NULL%1,               Han%Liao%NULL%1,               Sixuan%Long%NULL%1,               Ning%Zhou%NULL%0,               Pa%Wu%NULL%1]'
$s10 = '[Guangming%Ye%NULL%0,               Zhenyu%Pan%NULL%0,               Yunbao%Pan%NULL%0,               Qiaoling%Deng%NULL%0,               Liangjun%Chen%NULL%0,               Jin%Li%NULL%0,               Yirong%Li%NULL%0,               Xinghuan%Wang%NULL%0]'
$s11 = '[Jun%Yuan%NULL%1,               Zongqiu%Chen%NULL%1,               Chenghua%Gong%NULL%1,               Hui%Liu%NULL%1,               Baisheng%Li%NULL%1,               Kuibiao%Li%NULL%1,               Xi%Chen%NULL%4,               Conghui%Xu%NULL%1,               Qinlong%Jing%NULL%1,               Guocong%Liu%NULL%1,               Pengzhe%Qin%NULL%1,               Yufei%Liu%NULL%1,               Yi%Zhong%NULL%1,               Lijuan%Huang%NULL%1,               Bao-Ping%Zhu%NULL%1,               Zhicong%Yang%baopingzhu@yahoo.com%1]'
$s12 = '[Kenneth I.%Zheng%NULL%0,               Xiao-Bo%Wang%NULL%1,               Xiang-Hong%Jin%NULL%1,               Wen-Yue%Liu%NULL%1,               Feng%Gao%NULL%1,               Yong-Ping%Chen%NULL%1,               Ming-Hua%Zheng%zhengmh@wmu.edu.cn%0]'

$ws.Range("E2").Value = $s0
$ws.Range("E3").Value = $s1
$ws.Range("E4").Value = $s2
$ws.Range("E5").Value = $s3
$ws.Range("E6").Value = $s4
$ws.Range("E7").Value = $s5
$ws.Range("E8").Value = $s6
$ws.Range("E9").Value = $s7
$ws.Range("E10").Value = $s8
$ws.Range("E11").Value = $s9
$ws.Range("E12").Value = $s7
$ws.Range("E13").Value = $s10
$ws.Range("E14").Value = $s11
$ws.Range("E15").Value = $s12
